$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsredovisning")
$tbl = $ws.ListObjects.Item("Tabell3")

# Grow the table (and the underlying sheet range) from A1:C10 to A1:C13 so
# the totals row moves from row 10 down to row 13, leaving two fresh data
# rows (10 & 11) plus one blank spacer row (12) above the new totals row.
$tbl.Resize($ws.Range("A1:C13"))

# New row: 2024-01-23, 1.5 lessons, "Radera aktivitet"
$ws.Cells.Item(10, 1).Value = 45314
$ws.Cells.Item(10, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(10, 2).Value = 1.5
$ws.Cells.Item(10, 3).Value = "Radera aktivitet"

# New row: 2024-01-23, 2 lessons, "Hämta uppgifter (sida)"
$ws.Cells.Item(11, 1).Value = 45314
$ws.Cells.Item(11, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "Hämta uppgifter (sida)"

# Re-write the totals row, now at row 13 (row 12 stays empty).
$ws.Cells.Item(13, 1).Value = "Summa"
$ws.Cells.Item(13, 2).Formula = "=SUBTOTAL(109,Tabell3[Antal lekt])"

$ws.Activate()
$ws.Range("A12").Select()
